# Update "想去人数" (want-to-go count) figures on the 展览 and 全部类型 sheets.

$wb = $excel.ActiveWorkbook

# Sheet "展览": rows 2-4, column F
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 1295
$wsExpo.Range("F3").Value = 2806
$wsExpo.Range("F4").Value = 252

# Sheet "全部类型": rows 3, 4, 6, column F
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 1295
$wsAll.Range("F4").Value = 2806
$wsAll.Range("F6").Value = 252
